$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temp "text" cell (outside the used range) used to force numeric-looking
# strings (e.g. "1.003") into the Price column as literal text via a
# values-only paste, so the destination cell keeps its original (default)
# style instead of picking up a Text number format.
$tmp = $ws.Range("Z1")
$tmp.NumberFormat = "@"

$ws.Range("D2").Value = "27.554.32"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.847.05"
$ws.Range("E3").Value = "  -2.24%  "
$tmp.Value = "1.003"
$tmp.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -1.39%  "
$tmp.Value = "332.95"
$tmp.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.84%  "
$tmp.Value = "1.003"
$tmp.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -1.24%  "
$tmp.Value = "0.4636"
$tmp.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -1.42%  "
$tmp.Value = "0.3855"
$tmp.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -2.27%  "
$tmp.Value = "46.14"
$tmp.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -1.47%  "
$tmp.Value = "0.07920"
$tmp.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -1.18%  "
$tmp.Value = "0.9949"
$tmp.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -2.19%  "
$tmp.Value = "21.50"
$tmp.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "1.843.26"
$ws.Range("E13").Value = "  -2.48%  "
$tmp.Value = "5.921"
$tmp.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -1.09%  "
$tmp.Value = "7.108"
$tmp.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -0.86%  "
$tmp.Value = "1.003"
$tmp.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -1.44%  "
$tmp.Value = "88.87"
$tmp.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +1.04%  "
$tmp.Value = "0.06636"
$tmp.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -2.28%  "
$tmp.Value = "0.00001036"
$tmp.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.48%  "
$tmp.Value = "17.08"
$tmp.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.77%  "
$tmp.Value = "1.003"
$tmp.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "27.550.01"
$ws.Range("E22").Value = "  -1.57%  "
$tmp.Value = "5.383"
$tmp.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -2.30%  "
$tmp.Value = "10.92"
$tmp.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.64%  "
$tmp.Value = "2.298"
$tmp.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -2.80%  "
$tmp.Value = "158.14"
$tmp.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.69%  "
$tmp.Value = "19.53"
$tmp.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.55%  "
$tmp.Value = "2.102"
$tmp.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -0.01%  "
$tmp.Value = "5.401"
$tmp.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -1.73%  "
$tmp.Value = "119.79"
$tmp.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -1.43%  "
$tmp.Value = "0.9762"
$tmp.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +1.07%  "
$tmp.Value = "0.09408"
$tmp.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -1.66%  "
$tmp.Value = "3.581"
$tmp.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.92%  "
$tmp.Value = "5.285"
$tmp.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.42%  "
$tmp.Value = "1.347"
$tmp.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.07%  "
$tmp.Value = "0.06015"
$tmp.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -1.90%  "
$tmp.Value = "0.02224"
$tmp.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -1.25%  "
$tmp.Value = "8.291"
$tmp.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +0.86%  "
$tmp.Value = "1.181"
$tmp.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -3.00%  "
$tmp.Value = "0.5892"
$tmp.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -1.26%  "
$tmp.Value = "0.1861"
$tmp.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.23%  "
$tmp.Value = "10.30"
$tmp.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -0.40%  "
$tmp.Value = "1.244"
$tmp.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.70%  "
$tmp.Value = "0.5581"
$tmp.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -2.12%  "
$tmp.Value = "12.14"
$tmp.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.42%  "
$tmp.Value = "1.899"
$tmp.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -2.51%  "
$tmp.Value = "0.06679"
$tmp.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -2.80%  "
$tmp.Value = "110.87"
$tmp.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("E49").Value = "  -1.64%  "
$tmp.Value = "1.002"
$tmp.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -1.40%  "
$tmp.Value = "70.02"
$tmp.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.60%  "

# Clean up the scratch cell so it does not leave stray data behind.
$tmp.Clear()
$excel.CutCopyMode = $false
